$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price/Volume columns as Text so numeric-looking strings
# (e.g. "506.29", "2.292.57") are not silently coerced into numbers by Excel.
$priceVolRange = $ws.Range("D2:E50")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "54.847.12"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "2.292.57"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "506.29"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("D6").Value = "129.57"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "2.315.35"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").Value = "0.0979"
$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("E11").Value = "  +1.60%  "

$ws.Range("D12").Value = "5.06"
$ws.Range("E12").Value = "  +7.05%  "

$ws.Range("D13").Value = "0.341"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").Value = "23.77"
$ws.Range("E14").Value = "  +3.67%  "

$ws.Range("D15").Value = "2.699.57"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").Value = "54.849.06"
$ws.Range("E16").Value = "  +0.89%  "

$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("D18").Value = "2.293.46"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "10.74"
$ws.Range("E19").Value = "  +4.48%  "

$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "310.23"
$ws.Range("E21").Value = "  +2.04%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  +3.95%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "60.13"
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("D25").Value = "0.993"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").Value = "7.51"
$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("D28").Value = "173.22"
$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").Value = "6.13"
$ws.Range("E29").Value = "  +2.57%  "

$ws.Range("D30").Value = "0.0₃0708"
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "18.07"
$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.25%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "0.920"
$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("D38").Value = "3.88"
$ws.Range("E38").Value = "  +2.55%  "

$ws.Range("D39").Value = "36.76"
$ws.Range("E39").Value = "  +1.90%  "

$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("D41").Value = "0.377"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").Value = "135.49"
$ws.Range("E42").Value = "  +7.58%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "5.12"
$ws.Range("E43").Value = "  +3.27%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.42"
$ws.Range("E44").Value = "  +0.93%  "

$ws.Range("D45").Value = "257.71"
$ws.Range("E45").Value = "  +6.11%  "

$ws.Range("D46").Value = "0.0504"
$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("D47").Value = "0.0911"
$ws.Range("E47").Value = "  +1.66%  "

$ws.Range("D48").Value = "0.551"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("E50").Value = "  +1.72%  "

# Restore the default (Normal) style so cell formatting matches the original
# workbook; only the underlying text values should differ.
$priceVolRange.Style = "Normal"